$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '58.705.32'
$ws.Range("E2").Value = '  -3.86%  '
$ws.Range("D3").Value = '3.204.72'
$ws.Range("E3").Value = '  -4.80%  '
$ws.Range("D4").Value = '''0.999'
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").Value = '''533.77'
$ws.Range("E5").Value = '  -6.00%  '
$ws.Range("D6").Value = '''134.67'
$ws.Range("E6").Value = '  -9.54%  '
$ws.Range("E7").Value = '  -0.11%  '
$ws.Range("D8").Value = '3.202.49'
$ws.Range("E8").Value = '  -4.91%  '
$ws.Range("D9").Value = '''0.456'
$ws.Range("E9").Value = '  -5.02%  '
$ws.Range("D10").Value = '''7.54'
$ws.Range("E10").Value = '  -5.58%  '
$ws.Range("E11").Value = '  -7.01%  '
$ws.Range("D12").Value = '''0.391'
$ws.Range("E12").Value = '  -5.85%  '
$ws.Range("D13").Value = '3.755.62'
$ws.Range("E13").Value = '  -4.80%  '
$ws.Range("E14").Value = '  -0.67%  '
$ws.Range("D15").Value = '''25.80'
$ws.Range("E15").Value = '  -7.99%  '
$ws.Range("D16").Value = '3.208.80'
$ws.Range("E16").Value = '  -4.50%  '
$ws.Range("E17").Value = '  -7.43%  '
$ws.Range("D18").Value = '58.738.60'
$ws.Range("E18").Value = '  -3.88%  '
$ws.Range("D19").Value = '''5.90'
$ws.Range("E19").Value = '  -7.17%  '
$ws.Range("D20").Value = '''13.22'
$ws.Range("E20").Value = '  -8.62%  '
$ws.Range("D21").Value = '''8.19'
$ws.Range("E21").Value = '  -7.77%  '
$ws.Range("D22").Value = '''360.19'
$ws.Range("E22").Value = '  -4.24%  '
$ws.Range("E23").Value = '  -0.14%  '
$ws.Range("D24").Value = '''69.94'
$ws.Range("E24").Value = '  -7.22%  '
$ws.Range("D25").Value = '''0.517'
$ws.Range("E25").Value = '  -7.63%  '
$ws.Range("D26").Value = '3.346.09'
$ws.Range("E26").Value = '  -4.53%  '
$ws.Range("D27").Value = '''0.169'
$ws.Range("E27").Value = '  -3.85%  '
$ws.Range("D28").Value = '0.0₃0956'
$ws.Range("E28").Value = '  -11.79%  '
$ws.Range("D29").Value = '''1.00'
$ws.Range("E29").Value = '  +0.46%  '
$ws.Range("D30").Value = '''7.07'
$ws.Range("E30").Value = '  -5.08%  '
$ws.Range("D32").Value = '''1.92'
$ws.Range("E32").Value = '  -7.94%  '
$ws.Range("D33").Value = '''7.03'
$ws.Range("E33").Value = '  -8.86%  '
$ws.Range("D34").Value = '''21.67'
$ws.Range("E34").Value = '  -5.46%  '
$ws.Range("D35").Value = '''1.19'
$ws.Range("E35").Value = '  -7.89%  '
$ws.Range("B36").Value = 'Monero'
$ws.Range("C36").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D36").Value = '''161.52'
$ws.Range("E36").Value = '  -5.09%  '
$ws.Range("B37").Value = 'NEARProtocol'
$ws.Range("C37").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D37").Value = '''4.89'
$ws.Range("E37").Value = '  -9.20%  '
$ws.Range("D38").Value = '''6.34'
$ws.Range("E38").Value = '  -6.83%  '
$ws.Range("E39").Value = '  -8.02%  '
$ws.Range("D40").Value = '''25.79'
$ws.Range("E40").Value = '  -11.88%  '
$ws.Range("D41").Value = '''0.0703'
$ws.Range("E41").Value = '  -6.79%  '
$ws.Range("D42").Value = '3.232.83'
$ws.Range("E42").Value = '  -5.07%  '
$ws.Range("D43").Value = '''40.92'
$ws.Range("E43").Value = '  -3.54%  '
$ws.Range("D44").Value = '''0.714'
$ws.Range("E44").Value = '  -6.24%  '
$ws.Range("D45").Value = '''1.09'
$ws.Range("E45").Value = '  -4.65%  '
$ws.Range("D46").Value = '''4.00'
$ws.Range("E46").Value = '  -6.93%  '
$ws.Range("D47").Value = '''1.49'
$ws.Range("E47").Value = '  -6.98%  '
$ws.Range("D48").Value = '''0.998'
$ws.Range("E48").Value = '  -0.24%  '
$ws.Range("D49").Value = '2.294.57'
$ws.Range("D50").Value = '''6.26'
$ws.Range("E50").Value = '  -6.08%  '
$ws.Range("D51").Value = '''20.71'
$ws.Range("E51").Value = '  -8.22%  '
